$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = "$20.36"
    3  = "$25.91"
    4  = "$30.67"
    5  = "$43.09"
    6  = "$56.58"
    7  = "$70.06"
    8  = "$82.93"
    9  = "$236.93"
    10 = "$16.92"
    11 = "$21.41"
    12 = "$25.91"
    13 = "$39.66"
    14 = "$53.14"
    15 = "$70.06"
    16 = "$13.48"
    17 = "$16.92"
    18 = "$20.36"
    19 = "$31.73"
    20 = "$43.09"
    21 = "$52.08"
    22 = "$40.71"
    23 = "$52.08"
    24 = "$61.07"
    25 = "$85.92"
    26 = "$113.15"
    27 = "$113.15"
    28 = "$160.21"
    29 = "$327.30"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"

    $ws.Range("F$row").Value = "FAIL"
}
